$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# --- Crime statistics table updates (rows 14-31) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("F14").NumberFormat = "general"
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 7
$ws.Range("K14").Value = -28.571428571428
$ws.Range("N14").Value = -76.190476190476
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -40
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -38.888888888888
$ws.Range("L15").Value = -21.428571428571
$ws.Range("M15").Value = 37.5
$ws.Range("N15").Value = -60.714285714285
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 180
$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = 29.729729729729
$ws.Range("I16").Value = 151
$ws.Range("J16").Value = 180
$ws.Range("K16").Value = -16.111111111111
$ws.Range("L16").Value = -20.526315789473
$ws.Range("M16").Value = 4.137931034482
$ws.Range("N16").Value = -77.496274217585
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 23.529411764705
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = 2.597402597402
$ws.Range("I17").Value = 315
$ws.Range("J17").Value = 343
$ws.Range("K17").Value = -8.163265306122
$ws.Range("L17").Value = -13.698630136986
$ws.Range("M17").Value = 73.076923076923
$ws.Range("N17").Value = -17.105263157894
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 5.555555555555
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = -19
$ws.Range("L18").Value = -55.248618784530
$ws.Range("M18").Value = -3.571428571428
$ws.Range("N18").Value = -88.311688311688
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -43.243243243243
$ws.Range("I19").Value = 246
$ws.Range("J19").Value = 287
$ws.Range("K19").Value = -14.285714285714
$ws.Range("L19").Value = -1.992031872509
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 5.128205128205
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = -11.764705882352
$ws.Range("I20").Value = 115
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = 7.476635514018
$ws.Range("L20").Value = -31.952662721893
$ws.Range("M20").Value = 101.754385964912
$ws.Range("N20").Value = -78.504672897196
$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 49
$ws.Range("E21").Value = 38.775510204081
$ws.Range("F21").Value = 221
$ws.Range("G21").Value = 248
$ws.Range("H21").Value = -10.887096774193
$ws.Range("I21").Value = 924
$ws.Range("J21").Value = 1042
$ws.Range("K21").Value = -11.324376199616
$ws.Range("L21").Value = -21.428571428571
$ws.Range("M21").Value = 52.224052718286
$ws.Range("N21").Value = -63.934426229508
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -29.166666666666
$ws.Range("L22").Value = 6.25
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 19
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = -36.666666666666
$ws.Range("L23").Value = -29.629629629629
$ws.Range("M23").Value = -13.636363636363
$ws.Range("C24").Value = 39
$ws.Range("E24").Value = -15.217391304347
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 169
$ws.Range("H24").Value = -10.059171597633
$ws.Range("I24").Value = 645
$ws.Range("J24").Value = 679
$ws.Range("K24").Value = -5.007363770250
$ws.Range("L24").Value = -4.160475482912
$ws.Range("M24").Value = 51.764705882352
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -65.384615384615
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 91
$ws.Range("H25").Value = -52.747252747252
$ws.Range("I25").Value = 234
$ws.Range("J25").Value = 366
$ws.Range("K25").Value = -36.065573770491
$ws.Range("L25").Value = -31.378299120234
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 27
$ws.Range("E26").Value = -25.925925925925
$ws.Range("F26").Value = 78
$ws.Range("G26").Value = 106
$ws.Range("H26").Value = -26.415094339622
$ws.Range("I26").Value = 379
$ws.Range("J26").Value = 425
$ws.Range("K26").Value = -10.823529411764
$ws.Range("L26").Value = -4.292929292929
$ws.Range("M26").Value = 3.551912568306
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 28
$ws.Range("K27").Value = -35.714285714285
$ws.Range("L27").Value = -28
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 48
$ws.Range("J28").Value = 46
$ws.Range("K28").Value = 4.347826086956
$ws.Range("L28").Value = -4
$ws.Range("D29").Value = 2
$ws.Range("G29").Value = 9
$ws.Range("H29").Value = -77.777777777777
$ws.Range("J29").Value = 28
$ws.Range("K29").Value = -67.857142857142
$ws.Range("N29").Value = -83.333333333333
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -83.333333333333
$ws.Range("J30").Value = 18
$ws.Range("K30").Value = -66.666666666666
$ws.Range("N30").Value = -86.666666666666
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("D31").NumberFormat = "general"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("E31").NumberFormat = "general"
$ws.Range("G31").Value = 1
